$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 164, pushing existing rows 164:174 down to 165:175.
$ws.Rows.Item(164).Insert()

# Populate the newly inserted row 164 with the new data record.
$ws.Cells.Item(164, 1).Value = 10
$ws.Cells.Item(164, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(164, 3).Value = "La Araucanía"
$ws.Cells.Item(164, 4).Value = 44461
$ws.Cells.Item(164, 5).Value = 9
$ws.Cells.Item(164, 6).Value = 100112017
$ws.Cells.Item(164, 7).Value = "Apio"
$ws.Cells.Item(164, 8).Value = "Americana (o)"
$ws.Cells.Item(164, 9).Value = "Primera"
$ws.Cells.Item(164, 10).Value = 40
$ws.Cells.Item(164, 11).Value = 10000
$ws.Cells.Item(164, 12).Value = 10000
$ws.Cells.Item(164, 13).Value = 10000
$ws.Cells.Item(164, 14).Value = '$/docena de matas'
$ws.Cells.Item(164, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(164, 16).Value = 1667
$ws.Cells.Item(164, 17).Value = 6
$ws.Cells.Item(164, 18).Value = "Hortaliza"
